$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 1 de Abril de 2020 a las 21:55'

$ws.Range("B4").Value = 209105
$ws.Range("C4").Value = 20575
$ws.Range("D4").Value = 8805
$ws.Range("E4").Value = 195642
$ws.Range("F4").Value = 4946
$ws.Range("G4").Value = 605
$ws.Range("H4").Value = 4658

$ws.Range("B6").Value = 102179
$ws.Range("C6").Value = 6256
$ws.Range("E6").Value = 70401
$ws.Range("G6").Value = 667
$ws.Range("H6").Value = 9131

$ws.Range("B8").Value = 77558
$ws.Range("C8").Value = 5750
$ws.Range("E8").Value = 57967
$ws.Range("F8").Value = 3408
$ws.Range("G8").Value = 116
$ws.Range("H8").Value = 891

$ws.Range("B12").Value = 17768
$ws.Range("C12").Value = 1163
$ws.Range("E12").Value = 14313
$ws.Range("G12").Value = 55
$ws.Range("H12").Value = 488

$ws.Range("B16").Value = 10668
$ws.Range("C16").Value = 488
$ws.Range("E16").Value = 9086

$ws.Range("A32").Value = 'Polonia'
$ws.Range("B32").Value = 2554
$ws.Range("C32").Value = 243
$ws.Range("D32").Value = 7
$ws.Range("E32").Value = 2504
$ws.Range("F32").Value = 50
$ws.Range("H32").Value = 43

$ws.Range("A33").Value = 'Rumania'
$ws.Range("B33").Value = 2460
$ws.Range("C33").Value = 215
$ws.Range("D33").Value = 252
$ws.Range("E33").Value = 2116
$ws.Range("F33").Value = 57
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 92

$ws.Range("E45").Value = 891
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 38

$ws.Range("F50").Value = 47

$ws.Range("A59").Value = 'Ucrania'
$ws.Range("B59").Value = 794
$ws.Range("C59").Value = 149
$ws.Range("D59").Value = 13
$ws.Range("E59").Value = 761
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 20

$ws.Range("A60").Value = 'Estonia'
$ws.Range("B60").Value = 779
$ws.Range("C60").Value = 34
$ws.Range("D60").Value = 33
$ws.Range("E60").Value = 741
$ws.Range("F60").Value = 15
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 5

$ws.Range("A61").Value = 'Egipto'
$ws.Range("B61").Value = 779
$ws.Range("C61").Value = 69
$ws.Range("D61").Value = 179
$ws.Range("E61").Value = 548
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 52

$ws.Range("A62").Value = 'Hong Kong'
$ws.Range("B62").Value = 765
$ws.Range("C62").Value = 50
$ws.Range("D62").Value = 147
$ws.Range("E62").Value = 614
$ws.Range("F62").Value = 5
$ws.Range("H62").Value = 4

$ws.Range("A63").Value = 'Irak'
$ws.Range("B63").Value = 728
$ws.Range("C63").Value = 34
$ws.Range("D63").Value = 182
$ws.Range("E63").Value = 494
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 52

$ws.Range("A64").Value = 'Crucero'
$ws.Range("B64").Value = 712
$ws.Range("D64").Value = 603
$ws.Range("E64").Value = 98
$ws.Range("F64").Value = 15
$ws.Range("H64").Value = 11

$ws.Range("A65").Value = 'Nueva Zelanda'
$ws.Range("B65").Value = 708
$ws.Range("C65").Value = 61
$ws.Range("D65").Value = 83
$ws.Range("E65").Value = 624
$ws.Range("F65").Value = 2
$ws.Range("H65").Value = 1

$ws.Range("A102").Value = 'Nigeria'
$ws.Range("B102").Value = 174
$ws.Range("C102").Value = 39
$ws.Range("D102").Value = 9
$ws.Range("E102").Value = 163
$ws.Range("F102").Value = 0
$ws.Range("H102").Value = 2

$ws.Range("A103").Value = 'Islas Feroe'
$ws.Range("B103").Value = 173
$ws.Range("C103").Value = 4
$ws.Range("D103").Value = 75
$ws.Range("E103").Value = 98
$ws.Range("F103").Value = 1
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 0

$ws.Range("A104").Value = 'Honduras'
$ws.Range("B104").Value = 172
$ws.Range("C104").Value = 31
$ws.Range("D104").Value = 3
$ws.Range("E104").Value = 159
$ws.Range("F104").Value = 4
$ws.Range("G104").Value = 3
$ws.Range("H104").Value = 10

$ws.Range("A105").Value = 'Bielorrusia'
$ws.Range("B105").Value = 163
$ws.Range("C105").Value = 11
$ws.Range("D105").Value = 53
$ws.Range("E105").Value = 108
$ws.Range("F105").Value = 2
$ws.Range("H105").Value = 2

$ws.Range("A106").Value = 'Mauricio'
$ws.Range("B106").Value = 161
$ws.Range("C106").Value = 18
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 155
$ws.Range("F106").Value = 1
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 6

$ws.Range("A158").Value = 'Birmania'
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 0
$ws.Range("H158").Value = 1

$ws.Range("A160").Value = 'Bahamas'
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 1
$ws.Range("H160").Value = 0

$ws.Range("A163").Value = 'Mongolia'
$ws.Range("C163").Value = 2

$ws.Range("A164").Value = 'Namibia'
$ws.Range("C164").Value = 3

$ws.Range("A169").Value = 'Mozambique'
$ws.Range("C169").Value = 2

$ws.Range("A171").Value = 'Laos'
$ws.Range("C171").Value = 1

$ws.Range("A172").Value = 'Seychelles'
$ws.Range("C172").Value = 0

$ws.Range("A186").Value = 'Islas Turcas y Caicos'
$ws.Range("C186").Value = 1

$ws.Range("A187").Value = 'Santa Sede'
$ws.Range("C187").Value = 0

$ws.Range("A188").Value = 'Liberia'
$ws.Range("C188").Value = 3

$ws.Range("A192").Value = 'Fiyi'

$ws.Range("A193").Value = 'Montserrat'
